$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.001.98"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.726.42"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.55"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.40%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4851"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.37%  "
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.23"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07231"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("E11").Value = "  -3.04%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.92"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.867"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.738.32"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.803"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.93%  "
$ws.Range("E17").Value = "  -5.82%  "
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06400"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.56"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.704"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.067.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("E25").Value = "  -3.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.89"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.923.90"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.066"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.55"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.036"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09316"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.636"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.370"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05902"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("E36").Value = "  -4.02%  "
$ws.Range("E37").Value = "  +5.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.1993"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.91"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.737"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5968"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.117"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.475"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.72"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.581"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5602"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.23"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.843"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.30%  "
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06648"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.39%  "
